$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.147.23"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "1.855.78"
$ws.Range("E3").Value = "  -3.05%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.61"
$ws.Range("E5").Value = "  -2.68%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4695"
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2805"
$ws.Range("E8").Value = "  -1.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06536"
$ws.Range("E9").Value = "  -2.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.92"
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.65"
$ws.Range("E12").Value = "  -6.09%  "
$ws.Range("D13").Value = "1.854.16"
$ws.Range("E13").Value = "  -3.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.077"
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6651"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "282.34"
$ws.Range("E16").Value = "  +2.57%  "
$ws.Range("D17").Value = "30.181.13"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.473"
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.56"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").Value = "2.103.20"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007216"
$ws.Range("E22").Value = "  -3.56%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.117"
$ws.Range("E24").Value = "  -3.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.74"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.272"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.916"
$ws.Range("E28").Value = "  -7.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.343"
$ws.Range("E29").Value = "  -2.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09592"
$ws.Range("E30").Value = "  -3.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.410"
$ws.Range("E31").Value = "  -3.76%  "
$ws.Range("E32").Value = "  -3.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.088"
$ws.Range("E33").Value = "  -3.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04668"
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.090"
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6927"
$ws.Range("E36").Value = "  -4.79%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.708"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01851"
$ws.Range("E39").Value = "  -3.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.278"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.507"
$ws.Range("E41").Value = "  -4.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.65"
$ws.Range("E42").Value = "  -2.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8584"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.941"
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.96"
$ws.Range("E46").Value = "  -2.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4144"
$ws.Range("E47").Value = "  -2.81%  "
$ws.Range("D48").Value = "1.023.38"
$ws.Range("E48").Value = "  +7.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.190"
$ws.Range("E49").Value = "  -2.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.890"
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.67"
$ws.Range("E51").Value = "  -2.87%  "